$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated cross-section (spectroscopic factor) distributions for bite 3
# modified run properties. Each row below corresponds to one energy level
# whose L, SPECTROSCOPIC_FACTOR, and/or ERROR values were recomputed.

$ws.Cells.Item(2, 4).Value = 0.1003833467679155
$ws.Cells.Item(2, 5).Value = 0.002633670046383253

$ws.Cells.Item(3, 3).Value = 3
$ws.Cells.Item(3, 4).Value = 0.01967941761905094
$ws.Cells.Item(3, 5).Value = 0.001346880293319076

$ws.Cells.Item(4, 3).Value = 1
$ws.Cells.Item(4, 4).Value = 0.05943802801823975
$ws.Cells.Item(4, 5).Value = 0.001053883381133756

$ws.Cells.Item(5, 3).Value = 2
$ws.Cells.Item(5, 4).Value = 0.01003317029389575
$ws.Cells.Item(5, 5).Value = 0.001693911867800581

$ws.Cells.Item(6, 2).Value = 2552.302019868145
$ws.Cells.Item(6, 3).Value = 1
$ws.Cells.Item(6, 4).Value = 0.01383371637882464
$ws.Cells.Item(6, 5).Value = 0.0007008793442481233

$ws.Cells.Item(7, 3).Value = 2
$ws.Cells.Item(7, 4).Value = 0.01187085410814768
$ws.Cells.Item(7, 5).Value = 0.001427097215187148

$ws.Cells.Item(8, 2).Value = 2576.595223088504
$ws.Cells.Item(8, 3).Value = 1
$ws.Cells.Item(8, 4).Value = 0.01749281742997389
$ws.Cells.Item(8, 5).Value = 0.000670410960262062

$ws.Cells.Item(9, 3).Value = 3
$ws.Cells.Item(9, 4).Value = 0.008630323237560395
$ws.Cells.Item(9, 5).Value = 0.00112858073106559

$ws.Cells.Item(10, 3).Value = 1
$ws.Cells.Item(10, 4).Value = 0.02695315010985294
$ws.Cells.Item(10, 5).Value = 0.001059372864211382

$ws.Cells.Item(11, 3).Value = 1
$ws.Cells.Item(11, 4).Value = 0.003118765279171087
$ws.Cells.Item(11, 5).Value = 0.0006307310414617373

$ws.Cells.Item(12, 3).Value = 1
$ws.Cells.Item(12, 4).Value = 0.01153263396511821
$ws.Cells.Item(12, 5).Value = 0.0006702065028555223

$ws.Cells.Item(13, 3).Value = 1
$ws.Cells.Item(13, 4).Value = 0.01477815843567994
$ws.Cells.Item(13, 5).Value = 0.0007604345812864912

$ws.Cells.Item(14, 3).Value = 1
$ws.Cells.Item(14, 4).Value = 0.007289291180016085
$ws.Cells.Item(14, 5).Value = 0.0006598198330066742

$ws.Cells.Item(15, 3).Value = 1
$ws.Cells.Item(15, 4).Value = 0.01800662804138555
$ws.Cells.Item(15, 5).Value = 0.000931406133057303

$ws.Cells.Item(16, 3).Value = 2
$ws.Cells.Item(16, 4).Value = 0.01534716269474826
$ws.Cells.Item(16, 5).Value = 0.0006263605004839691

$ws.Cells.Item(17, 3).Value = 1
$ws.Cells.Item(17, 4).Value = 0.0104679443723541
$ws.Cells.Item(17, 5).Value = 0.0005962120273638413

$ws.Cells.Item(18, 3).Value = 1
$ws.Cells.Item(18, 4).Value = 0.0278026913402598
$ws.Cells.Item(18, 5).Value = 0.001103103675507483

$ws.Cells.Item(19, 3).Value = 1
$ws.Cells.Item(19, 4).Value = 0.01406444463080108
$ws.Cells.Item(19, 5).Value = 0.001024149975157255

$ws.Cells.Item(20, 3).Value = 2
$ws.Cells.Item(20, 4).Value = 0.01644381189109301
$ws.Cells.Item(20, 5).Value = 0.001688147988819468

$ws.Cells.Item(21, 3).Value = 1
$ws.Cells.Item(21, 4).Value = 0.04967761623345822
$ws.Cells.Item(21, 5).Value = 0.001238525999353802

$ws.Cells.Item(22, 3).Value = 2
$ws.Cells.Item(22, 4).Value = 0.004770638388218259
$ws.Cells.Item(22, 5).Value = 0.0008673887978578587

$ws.Cells.Item(23, 3).Value = 2
$ws.Cells.Item(23, 4).Value = 0.001173930436090106
$ws.Cells.Item(23, 5).Value = 0.0006796439366837606

$ws.Cells.Item(24, 3).Value = 1
$ws.Cells.Item(24, 4).Value = 0.01786830746375652
$ws.Cells.Item(24, 5).Value = 0.0006846600141425562

$ws.Cells.Item(25, 3).Value = 1
$ws.Cells.Item(25, 4).Value = 0.01872632096582117
$ws.Cells.Item(25, 5).Value = 0.001673442769643433

$ws.Cells.Item(26, 4).Value = 0.0147537232516338
$ws.Cells.Item(26, 5).Value = 0.003902252889019084

$ws.Cells.Item(27, 3).Value = 2
$ws.Cells.Item(27, 4).Value = 0.006905271629601476
$ws.Cells.Item(27, 5).Value = 0.001038846174364854

$ws.Cells.Item(28, 4).Value = 0.006815869725090891
$ws.Cells.Item(28, 5).Value = 0.002039394090972078

$ws.Cells.Item(29, 4).Value = 0.03107694530033408
$ws.Cells.Item(29, 5).Value = 0.00204312097130224

$ws.Cells.Item(30, 4).Value = 0.04964234777733157
$ws.Cells.Item(30, 5).Value = 0.001996485725827465

$ws.Cells.Item(31, 3).Value = 4
$ws.Cells.Item(31, 4).Value = 0.007916632780802336
$ws.Cells.Item(31, 5).Value = 0.0008633424786485272

$ws.Cells.Item(32, 4).Value = 0.0159103685394359
$ws.Cells.Item(32, 5).Value = 0.002823683153756541

$ws.Cells.Item(33, 3).Value = 1
$ws.Cells.Item(33, 4).Value = 0.01788892509131872
$ws.Cells.Item(33, 5).Value = 0.001245416013572331

$ws.Cells.Item(34, 4).Value = 0.03463617967285389
$ws.Cells.Item(34, 5).Value = 0.002913837337557551

$ws.Cells.Item(35, 4).Value = 0.01025786663621027
$ws.Cells.Item(35, 5).Value = 0.001599344798118806

$ws.Cells.Item(36, 3).Value = 3
$ws.Cells.Item(36, 4).Value = 0.007476346500770412
$ws.Cells.Item(36, 5).Value = 0.001027284252014255

$ws.Cells.Item(37, 4).Value = 0.01237648045970882
$ws.Cells.Item(37, 5).Value = 0.00161674744743944

$ws.Cells.Item(38, 4).Value = 0.003973759564147854
$ws.Cells.Item(38, 5).Value = 0.002406643116314897

$ws.Cells.Item(39, 3).Value = 3
$ws.Cells.Item(39, 4).Value = 0.00385444869379247
$ws.Cells.Item(39, 5).Value = 0.0009496467796300289

$ws.Cells.Item(40, 3).Value = 2
$ws.Cells.Item(40, 4).Value = 0.008892275482153213
$ws.Cells.Item(40, 5).Value = 0.000994530810504012

$ws.Cells.Item(41, 4).Value = 0.02308051077533096
$ws.Cells.Item(41, 5).Value = 0.002778734712509133

$ws.Cells.Item(42, 3).Value = 2
$ws.Cells.Item(42, 4).Value = 0.003778260321143645
$ws.Cells.Item(42, 5).Value = 0.001104414555411244

$ws.Cells.Item(43, 3).Value = 2
$ws.Cells.Item(43, 4).Value = 0.002493016360014469
$ws.Cells.Item(43, 5).Value = 0.0002627792643487508

$ws.Cells.Item(44, 2).Value = 3259.461277156506
$ws.Cells.Item(44, 3).Value = 3
$ws.Cells.Item(44, 4).Value = 0.003520979003111298
$ws.Cells.Item(44, 5).Value = 0.001029209247063302

$ws.Cells.Item(45, 3).Value = 1
$ws.Cells.Item(45, 4).Value = 0.006690400924903099
$ws.Cells.Item(45, 5).Value = 0.0004779330702489417

$ws.Cells.Item(46, 3).Value = 3
$ws.Cells.Item(46, 4).Value = 0.007114685598151072
$ws.Cells.Item(46, 5).Value = 0.0008559020268452416

$ws.Cells.Item(47, 2).Value = 3322.739257221812
$ws.Cells.Item(47, 3).Value = 1
$ws.Cells.Item(47, 4).Value = 0.004031759145448133
$ws.Cells.Item(47, 5).Value = 0.000309780126409277

$ws.Cells.Item(48, 3).Value = 1
$ws.Cells.Item(48, 4).Value = 0.008720165543479061
$ws.Cells.Item(48, 5).Value = 0.000382519680734768

$ws.Cells.Item(49, 3).Value = 1
$ws.Cells.Item(49, 4).Value = 0.01098980593122331
$ws.Cells.Item(49, 5).Value = 0.000464825000048043

$ws.Cells.Item(50, 3).Value = 3
$ws.Cells.Item(50, 4).Value = 0.007028112719978081
$ws.Cells.Item(50, 5).Value = 0.001509742732439736

$ws.Cells.Item(51, 3).Value = 1
$ws.Cells.Item(51, 4).Value = 0.006788603562007698
$ws.Cells.Item(51, 5).Value = 0.0005157876535650174

$ws.Cells.Item(52, 3).Value = 1
$ws.Cells.Item(52, 4).Value = 0.01316921060467234
$ws.Cells.Item(52, 5).Value = 0.0004666798297619102

$ws.Cells.Item(53, 3).Value = 1
$ws.Cells.Item(53, 4).Value = 0.008655787788784879
$ws.Cells.Item(53, 5).Value = 0.0003797168811668873

$ws.Cells.Item(54, 3).Value = 3
$ws.Cells.Item(54, 4).Value = 0.006101726008472538
$ws.Cells.Item(54, 5).Value = 0.0009152589012708806

$ws.Cells.Item(55, 4).Value = 0.006343840521350375
$ws.Cells.Item(55, 5).Value = 0.001585960130337594

$ws.Cells.Item(56, 2).Value = 3485.001985506649
$ws.Cells.Item(56, 3).Value = 3
$ws.Cells.Item(56, 4).Value = 0.006843909527091864
$ws.Cells.Item(56, 5).Value = 0.0009058115550562761

$ws.Cells.Item(57, 3).Value = 3
$ws.Cells.Item(57, 4).Value = 0.01364478371902552
$ws.Cells.Item(57, 5).Value = 0.001049598747617347

$ws.Cells.Item(58, 3).Value = 1
$ws.Cells.Item(58, 4).Value = 0.005134959864739125
$ws.Cells.Item(58, 5).Value = 0.0006544222277116336

$ws.Cells.Item(59, 2).Value = 3536.043535477975
$ws.Cells.Item(59, 4).Value = 0.00612602588370704
$ws.Cells.Item(59, 5).Value = 0.00637606775651141

$ws.Cells.Item(60, 4).Value = 0.04215767530712843
$ws.Cells.Item(60, 5).Value = 0.00664987140261252

$ws.Cells.Item(61, 3).Value = 3
$ws.Cells.Item(61, 4).Value = 0.007645409937961853
$ws.Cells.Item(61, 5).Value = 0.001078198837404877

$ws.Cells.Item(62, 3).Value = 3
$ws.Cells.Item(62, 4).Value = 0.008299175502610256
$ws.Cells.Item(62, 5).Value = 0.0006840855412917909

$ws.Cells.Item(63, 4).Value = 0.03162934850620898
$ws.Cells.Item(63, 5).Value = 0.002231909672817166
